# Fix the "XML Tags" reference table:
#  - correct the capitalisation of the "Finding_Type" item name
#  - add a "Line_Number" row that reads from
#    AnalysisInfo/Unified/Context/FunctionDeclarationSourceLocation/line
#    (an Attribute), inserted right after the "File_Name" row
#  - the three Weakness_ID_* rows shift down by one row to make room

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("XML Tags")

# Grab the "plain border, no fill/alignment" formatting that row 7
# (the old Line_Number row) carries on its B/C cells, and stamp it onto
# C4 before that row's original content is overwritten below - this is
# the formatting the relocated "Attribute" cell ends up with.
$ws.Range("C7").Copy()
$ws.Range("C4").PasteSpecial(-4122)

# Row 4 becomes the relocated Line_Number / Attribute row.
$ws.Range("A4").Value = "Line_Number"
$ws.Range("B4").Value = "AnalysisInfo/Unified/Context/FunctionDeclarationSourceLocation/line"
$ws.Range("C4").Value = "Attribute"

# Item column: fix "Finding_type" -> "Finding_Type"
$ws.Range("A2").Value = "Finding_Type"

# Rows 5-7: the Weakness_ID_1..3 rows, each shifted down by one.
$ws.Range("A5").Value = "Weakness_ID_1"
$ws.Range("B5").Value = "ClassInfo/Kingdom"
$ws.Range("C5").Value = "Tag"

$ws.Range("A6").Value = "Weakness_ID_2"
$ws.Range("B6").Value = "ClassInfo/Type"
$ws.Range("C6").Value = "Tag"

$ws.Range("A7").Value = "Weakness_ID_3"
$ws.Range("B7").Value = "ClassInfo/Subtype"
$ws.Range("C7").Value = "Tag"

# The "XML Schema" column is now noticeably wider (longest entry is the
# new AnalysisInfo/.../line path), so the best-fit column widths grow.
$ws.Columns.Item(1).ColumnWidth = 15 - 0.8333333333333334
$ws.Columns.Item(2).ColumnWidth = 66 - 0.8333333333333334
$ws.Columns.Item(3).ColumnWidth = 15 - 0.8333333333333334

# Cursor/selection ends up on B16 after the edit.
$ws.Range("B16").Select()
